$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header updates ---
# Write N1 before L1 so new shared strings are appended in the order
# SlotPath100x100, then ItemMeshPathBoxMeshPath (matches target sharedStrings order).
$ws.Range("N1").Value = "SlotPath100x100"
$ws.Range("L1").Value = "ItemMeshPathBoxMeshPath"
$ws.Range("M1").Value = ""
$ws.Range("O1").Value = "ImagePath"
$ws.Range("K1").Value = "PathFile"

# --- RowName (col A) and ItemName (col D) values shift up the list by one ---
$ws.Range("A2").Value = "Name"
$ws.Range("A3").Value = "Won"
$ws.Range("A4").Value = "Gold"
$ws.Range("A5").Value = "Gold_Income"
$ws.Range("A6").Value = "Orange"
$ws.Range("A7").Value = "Apple"
$ws.Range("A8").Value = "Fish"
$ws.Range("A9").Value = "JuiceA"
$ws.Range("A10").Value = "MilkA"
$ws.Range("A11").Value = "MilkB"
$ws.Range("A12").Value = "Dumpling"
$ws.Range("A13").Value = "Hamburger"
$ws.Range("A14").Value = "SnackA"
$ws.Range("A15").Value = "SnackB"
$ws.Range("A16").Value = "SnackC"
$ws.Range("A17").Value = "IceCream"

$ws.Range("D2").Value = "Won"
$ws.Range("D3").Value = "Gold"
$ws.Range("D4").Value = "Gold_Income"
$ws.Range("D5").Value = "Orange"
$ws.Range("D6").Value = "Apple"
$ws.Range("D7").Value = "Fish"
$ws.Range("D8").Value = "JuiceA"
$ws.Range("D9").Value = "MilkA"
$ws.Range("D10").Value = "MilkB"
$ws.Range("D11").Value = "Dumpling"
$ws.Range("D12").Value = "Hamburger"
$ws.Range("D13").Value = "SnackA"
$ws.Range("D14").Value = "SnackB"
$ws.Range("D15").Value = "SnackC"
$ws.Range("D16").Value = "IceCream"

# --- K/L/M/N/O data columns (rows 2-16) ---
# Rows 2-4: all zeros
$ws.Range("K2:M4").Value = 0
$ws.Range("O2:O4").Value = 0

# Rows 5-16: K=8001, M=6001, N=7001 constant; L and O vary per row
$ws.Range("K5:K16").Value = 8001
$ws.Range("M5:M16").Value = 6001
$ws.Range("N5:N16").Value = 7001

$ws.Range("L5").Value = 5006
$ws.Range("L6").Value = 5004
$ws.Range("L7").Value = 5005
$ws.Range("L8").Value = 5008
$ws.Range("L9").Value = 5009
$ws.Range("L10").Value = 5007
$ws.Range("L11").Value = 5001
$ws.Range("L12").Value = 5002
$ws.Range("L13").Value = 5010
$ws.Range("L14").Value = 5011
$ws.Range("L15").Value = 5012
$ws.Range("L16").Value = 5003

# O5:O16 unchanged (101-112) already in file - leave as-is

# --- Row 17: A17 value, O17 blank (already blank) ---
$ws.Range("A17").Value = "IceCream"

# --- Column width cleanup: remove custom widths for columns K, L, M ---
$ws.Range("K1").EntireColumn.ColumnWidth = 8.43
$ws.Range("L1").EntireColumn.ColumnWidth = 8.43
$ws.Range("M1").EntireColumn.ColumnWidth = 8.43

# --- Selection change ---
$ws.Range("S6").Select()
